$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Delete the entire "is_global_admin" column (AH), shifting all columns to its right left by one.
$ws.Range("AH1").EntireColumn.Delete() | Out-Null

# Reflect the final selection state left behind by the author (cell AK2, which now
# holds the former "last_signin_ip" header column after the shift).
$ws.Range("AK2").Select() | Out-Null
